$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for new column F
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("F1").Value = "time_taken"

# Timestamps for each data row (row 2 through row 31)
$timestamps = @(
    "2021-10-05 13:39:59.522723",
    "2021-10-05 13:39:59.522736",
    "2021-10-05 13:39:59.522740",
    "2021-10-05 13:39:59.522743",
    "2021-10-05 13:39:59.522747",
    "2021-10-05 13:39:59.522750",
    "2021-10-05 13:39:59.522753",
    "2021-10-05 13:39:59.522756",
    "2021-10-05 13:39:59.522759",
    "2021-10-05 13:39:59.522763",
    "2021-10-05 13:39:59.522766",
    "2021-10-05 13:39:59.522768",
    "2021-10-05 13:39:59.522771",
    "2021-10-05 13:39:59.522775",
    "2021-10-05 13:39:59.522778",
    "2021-10-05 13:39:59.522781",
    "2021-10-05 13:39:59.522784",
    "2021-10-05 13:39:59.522787",
    "2021-10-05 13:39:59.522790",
    "2021-10-05 13:39:59.522793",
    "2021-10-05 13:39:59.522796",
    "2021-10-05 13:39:59.522799",
    "2021-10-05 13:39:59.522802",
    "2021-10-05 13:39:59.522805",
    "2021-10-05 13:39:59.522809",
    "2021-10-05 13:39:59.522812",
    "2021-10-05 13:39:59.522815",
    "2021-10-05 13:39:59.522818",
    "2021-10-05 13:39:59.522821",
    "2021-10-05 13:39:59.522825"
)

for ($i = 0; $i -lt $timestamps.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $timestamps[$i]
}
